$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 24 ---
# (pushes the old blank row 24 -> 25, and the old "Total" row 25 -> 26)
$ws.Rows("24").Insert()

# --- Match formatting of the new row to the rest of the comparison table ---
# Column B uses the "label" style (like B6:B23); columns C:F use the "value" style (like C6:F23)
$ws.Range("B9").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C9:F9").Copy()
$ws.Range("C24:F24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the new row 24: "Sleeved Joint / Lug Joint" ---
$ws.Range("B24").Value2 = "Sleeved Joint / Lug Joint"
$ws.Range("C24").Formula = "=2*0.155"
$ws.Range("D24").Formula = "=2*0.107"
$ws.Range("E24").Formula = "=C24-D24"
$ws.Range("F24").Formula = "=(E24/C24)*100"

# --- Extend the conditional formatting that highlights negative "Reduction" values ---
$cf = $ws.Range("E6:E23").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("E6:E24"))

# --- Update the "Total" row (now row 26) so it sums through the new row ---
$ws.Range("C26").Formula = "=SUM(C6:C24)"
$ws.Range("D26").Formula = "=SUM(D6:D24)"
$ws.Range("E26").Formula = "=C26-D26"
$ws.Range("F26").Formula = "=(E26/C26)*100"

# --- Update the saved selection / active cell, as in the authored workbook ---
$null = $ws.Range("H25").Select()
